$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 338, shifting the existing rows 338-354 down to 339-355.
$ws.Rows.Item(338).Insert()

# Populate the newly inserted row 338 with the new weekly data point.
$ws.Cells.Item(338, 1).Value = 6
$ws.Cells.Item(338, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(338, 3).Value = "Metropolitana"
$ws.Cells.Item(338, 4).Value = 45147
$ws.Cells.Item(338, 5).Value = 13
$ws.Cells.Item(338, 6).Value = 100112029
$ws.Cells.Item(338, 7).Value = "Orégano"
$ws.Cells.Item(338, 8).Value = "Sin especificar"
$ws.Cells.Item(338, 9).Value = "Primera"
$ws.Cells.Item(338, 10).Value = 32
$ws.Cells.Item(338, 11).Value = 20000
$ws.Cells.Item(338, 12).Value = 20000
$ws.Cells.Item(338, 13).Value = 20000
$ws.Cells.Item(338, 14).Value = "`$/docena de atados"
$ws.Cells.Item(338, 15).Value = "Región Metropolitana"
$ws.Cells.Item(338, 16).Value = 6667
$ws.Cells.Item(338, 17).Value = 3
$ws.Cells.Item(338, 18).Value = "Hortaliza"
